$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# jumpers_map: new row documenting R77 (additional load for STEP-UP DC-DC LM2577)
$ws.Range("B24").Value = "R77"
$ws.Range("C24").Value = "Дополнительная нагрузка для STEP-UP DC-DC (LM2577)"
$ws.Range("D24").Value = "TDA2003V"
$ws.Range("E24").Value = "NC"

# "Назначение" column uses the wrapped/centered style already used by the
# other long-text cells in this table (e.g. C8, C23)
$ws.Range("C24").WrapText = $true

# Row grows to fit the wrapped text, matching row 23 right above it
$ws.Rows.Item(24).RowHeight = 30.75

# Leave the cursor where the author left it after typing the new row
$ws.Range("C28").Select() | Out-Null
